# Refresh of the CapitalIQ "in_process_ids" ticker list ("full data from 2000").
# The underlying CIQ add-in re-pulled the universe: Cardinal Health (CAH) has
# dropped off the list and Caesars Entertainment (CZR) has been added (sorted
# in alphabetically by company name, between "Cadence Design Systems, Inc."
# and "Campbell Soup Company").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new CZR row right after the CDNS row (row 83), pushing
#        CPB/COF down by one, and populate its ID / IQID / IQ Name columns.
$ws.Rows.Item(84).Insert()
$ws.Cells.Item(84, 1).Value = "CZR"
$ws.Cells.Item(84, 3).Value = "IQ3133890"
$ws.Cells.Item(84, 5).Value = "Caesars Entertainment, Inc."

# --- 2. Remove the old Cardinal Health (CAH) row, which is now the last row
#        of the table (row 87, after the insertion above).
$ws.Rows.Item(87).Delete()

# --- 3. Stamp the refreshed CIQ bookkeeping defined names (new workbook GUID
#        assigned by the CIQ add-in, and the new "names last refreshed" date).
$wb.Names.Item("CIQWBGuid").Value = '="cee513ec-c86f-4aaa-b3c4-9d98e69b96d0"'
$wb.Names.Item("IQ_NAMES_REVISION_DATE_").Value = "=44278.8607986111"
